# premier league - atualizacao de dados
# atualizando rodada 5 do sabado - 22 de setembro
#
# Updates corner-kick statistics for teams that played their 5th match
# (round 5, Saturday Sept 22). For each affected club: "partidas jogadas"
# (B) goes from 4 to 5, either "partidas como mandante" (C) or "partidas
# como visitante" (F) goes from 2 to 3, the related totals/averages
# (D/E or G/H) are updated, and the overall total (I) and average (J)
# are recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Aston Villa (row 3) - new home match
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 16
$ws.Range("E3").Value = 5.333333333333333
$ws.Range("I3").Value = 23
$ws.Range("J3").Value = 4.6

# Bournemouth (row 4) - new away match
$ws.Range("B4").Value = 5
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 6.333333333333333
$ws.Range("H4").Value = 19
$ws.Range("I4").Value = 33
$ws.Range("J4").Value = 6.6

# Brentford (row 5) - new away match
$ws.Range("B5").Value = 5
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3.666666666666667
$ws.Range("H5").Value = 11
$ws.Range("I5").Value = 17
$ws.Range("J5").Value = 3.4

# Chelsea (row 7) - new away match
$ws.Range("B7").Value = 5
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 4.333333333333333
$ws.Range("H7").Value = 13
$ws.Range("I7").Value = 21
$ws.Range("J7").Value = 4.2

# Crystal Palace (row 9) - new home match
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 12
$ws.Range("I9").Value = 23
$ws.Range("J9").Value = 4.6

# Everton (row 10) - new away match
$ws.Range("B10").Value = 5
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 2.666666666666667
$ws.Range("H10").Value = 8
$ws.Range("I10").Value = 17
$ws.Range("J10").Value = 3.4

# Fulham (row 11) - new home match
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 16
$ws.Range("E11").Value = 5.333333333333333
$ws.Range("I11").Value = 30

# Ipswich Town (row 12) - new away match
$ws.Range("B12").Value = 5
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 4.333333333333333
$ws.Range("H12").Value = 13
$ws.Range("I12").Value = 23
$ws.Range("J12").Value = 4.6

# Leicester City (row 13) - new home match
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 9
$ws.Range("E13").Value = 3
$ws.Range("I13").Value = 16
$ws.Range("J13").Value = 3.2

# Liverpool (row 14) - new home match
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 19
$ws.Range("E14").Value = 6.333333333333333
$ws.Range("I14").Value = 31
$ws.Range("J14").Value = 6.2

# Manchester United (row 15) - new away match
$ws.Range("B15").Value = 5
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 7.333333333333333
$ws.Range("H15").Value = 22
$ws.Range("I15").Value = 34
$ws.Range("J15").Value = 6.8

# Newcastle (row 16) - new away match
$ws.Range("B16").Value = 5
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 5.333333333333333
$ws.Range("J16").Value = 5.2

# Southampton (row 18) - new home match
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 5.2

# Tottenham (row 19) - new home match
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = 9.333333333333334
$ws.Range("I19").Value = 53
$ws.Range("J19").Value = 10.6

# West Ham (row 20) - new home match
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = 4.666666666666667
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 3.8
